$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "مازن احمد خيرت "
$ws.Range("B3").Value = "mazenahmeddd1233@gmail.com"
$ws.Range("C3").Value = "https://github.com/mohamedbelal42/portfolio3.git"

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:mazenahmeddd1233@gmail.com") | Out-Null
$ws.Range("B3").Style = "Hyperlink"

$ws.Range("C3").Select()
